$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new search-result rows to the existing table.
$ws.Range("A4").Value = "jackets"
$ws.Range("B4").Value = "0 results have been found."

$ws.Range("A5").Value = "top"
$ws.Range("B5").Value = "1 result has been found."

# Keep the active selection consistent with the new last row, like Excel would.
$ws.Range("A5").Select()
